# Updates cryptos list values (price/volume columns) per upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.922.07"
$ws.Range("E2").Value = "  +3.47%  "

$ws.Range("D3").Value = "2.643.50"
$ws.Range("E3").Value = "  +5.66%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'113.61"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +6.99%  "

$ws.Range("D6").Value = "'326.33"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.62%  "

$ws.Range("D7").Value = "'0.529"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.85%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("E9").Value = "  +3.22%  "

$ws.Range("D10").Value = "'40.91"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.93%  "

$ws.Range("D11").Value = "'20.07"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.05%  "

$ws.Range("E12").Value = "  +2.23%  "

$ws.Range("E13").Value = "  +0.63%  "

$ws.Range("E14").Value = "  +3.75%  "

$ws.Range("D15").Value = "3.061.92"
$ws.Range("E15").Value = "  +5.89%  "

$ws.Range("D16").Value = "2.662.29"
$ws.Range("E16").Value = "  +6.61%  "

$ws.Range("D17").Value = "'0.871"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.72%  "

$ws.Range("D18").Value = "49.862.63"
$ws.Range("E18").Value = "  +3.71%  "

$ws.Range("D19").Value = "'13.19"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.50%  "

$ws.Range("D20").Value = "'6.74"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.06%  "

$ws.Range("E21").Value = "  -3.13%  "

$ws.Range("E22").Value = "  +2.71%  "

$ws.Range("D23").Value = "'72.21"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.45%  "

$ws.Range("D24").Value = "'275.38"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.40%  "

$ws.Range("D26").Value = "'26.81"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.71%  "

$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("D28").Value = "'10.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.67%  "

$ws.Range("E29").Value = "  -0.90%  "

$ws.Range("D30").Value = "'36.41"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +5.02%  "

$ws.Range("E31").Value = "  +1.08%  "

$ws.Range("D32").Value = "'50.17"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.67%  "

$ws.Range("D33").Value = "'5.45"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.96%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0818"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.58%  "

$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").Value = "'19.47"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.56%  "

$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("D37").Value = "'5.03"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +9.44%  "

$ws.Range("D38").Value = "'2.06"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +6.11%  "

$ws.Range("E39").Value = "  +7.62%  "

$ws.Range("D40").Value = "'123.99"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.50%  "

$ws.Range("E41").Value = "  +1.77%  "

$ws.Range("E42").Value = "  +0.57%  "

$ws.Range("D43").Value = "'21.98"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.52%  "

$ws.Range("E44").Value = "  +4.57%  "

$ws.Range("D45").Value = "2.087.12"
$ws.Range("E45").Value = "  +4.17%  "

$ws.Range("D46").Value = "'3.32"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.53%  "

$ws.Range("D47").Value = "'2.26"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +13.19%  "

$ws.Range("E48").Value = "  +4.59%  "

$ws.Range("E49").Value = "  +2.19%  "

$ws.Range("D50").Value = "'5.38"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.62%  "

$ws.Range("D51").Value = "'59.71"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.53%  "
